$p = $ppt.ActivePresentation
$m = $p.SlideMaster
Write-Output "Before: $($m.Name)"
$m.Name = "Office Theme"
Write-Output "After: $($m.Name)"
